$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Updated daily COVID-19 figures for Valais (rows 215-251),
# matching the source daily-tracking values as of the new upload.
# Columns: C=nouveaux cas positifs, E=patients SI, F=patients intubes,
# G=patients hospitalises hors SI, L=deces hopital, M=deces extra-hospitaliers.
# B/H/J/K recompute automatically (shared TODAY()-driven formulas).
$dailyUpdates = @{
    215 = @{ G = 13 }
    216 = @{ G = 13 }
    217 = @{ G = 16 }
    218 = @{ G = 17 }
    219 = @{ G = 19 }
    220 = @{ G = 20 }
    221 = @{ G = 21 }
    222 = @{ G = 24 }
    223 = @{ G = 26 }
    224 = @{ G = 26 }
    225 = @{ G = 28 }
    226 = @{ G = 30 }
    227 = @{ G = 35 }
    228 = @{ G = 38 }
    229 = @{ G = 42 }
    230 = @{ G = 45 }
    231 = @{ G = 55 }
    232 = @{ G = 65 }
    233 = @{ C = 310; G = 70 }
    234 = @{ C = 397; G = 81 }
    235 = @{ C = 283; G = 83 }
    236 = @{ G = 83 }
    237 = @{ G = 76 }
    238 = @{ C = 559; G = 79 }
    239 = @{ C = 669; G = 88 }
    240 = @{ C = 678; G = 110 }
    241 = @{ C = 761; G = 124 }
    242 = @{ G = 143 }
    243 = @{ C = 376; G = 142 }
    244 = @{ C = 881; G = 160 }
    245 = @{ C = 886; G = 172 }
    246 = @{ C = 767; G = 173 }
    247 = @{ C = 785; E = 15; G = 184; L = 4; M = 5 }
    248 = @{ C = 793; E = 18; G = 197 }
    249 = @{ C = 454; E = 22; F = 7; G = 211; L = 0; M = 0 }
    250 = @{ C = 277; E = 21; F = 6; G = 228; L = 2; M = 0 }
    251 = @{ C = 10; E = 24; F = 9; G = 247; L = 0; M = 0 }
}

foreach ($row in $dailyUpdates.Keys) {
    $rowData = $dailyUpdates[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}

# Move the active selection on the frozen (bottom-right) pane to the
# cell where data entry finished.
$ws.Range("H256").Select()
